$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("AI3").Value = 19

# Row 4
$ws.Range("AQ4").Value = 1000

# Row 5
$ws.Range("G5").Value = 2.1
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 3.4
$ws.Range("J5").Value = 3
$ws.Range("L5").Value = 4.5
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 7
$ws.Range("Q5").Value = 1.98
$ws.Range("R5").Value = 1.88
$ws.Range("S5").Value = 2.6
$ws.Range("T5").Value = 1.48
$ws.Range("AD5").Value = 19
$ws.Range("AL5").Value = 17
$ws.Range("AR5").Value = 4.3
$ws.Range("AS5").Value = 1.21

# Row 8
$ws.Range("G8").Value = 2.38
$ws.Range("H8").Value = 3.2
$ws.Range("I8").Value = 2.9
$ws.Range("J8").Value = 3.25
$ws.Range("L8").Value = 4
$ws.Range("M8").Value = 1.11
$ws.Range("N8").Value = 6.5
$ws.Range("Y8").Value = 2.2
$ws.Range("Z8").Value = 1.62
$ws.Range("AA8").Value = 6
$ws.Range("AB8").Value = 10
$ws.Range("AC8").Value = 11
$ws.Range("AD8").Value = 23
$ws.Range("AL8").Value = 13
$ws.Range("AM8").Value = 12
$ws.Range("AN8").Value = 34
$ws.Range("AO8").Value = 29
